$wb = $excel.ActiveWorkbook

# Rename sheets: mars -> March, april -> April
$wsMarch = $wb.Worksheets.Item("mars")
$wsMarch.Name = "March"

$wsApril = $wb.Worksheets.Item("april")
$wsApril.Name = "April"

# A new budget entry was added for an existing category (Transportation) and
# month (April). It is inserted before the existing "Monthly total" row,
# pushing that row down, and the total is updated to reflect the new entry.
$wsApril.Rows.Item(10).Insert()

$wsApril.Range("A10").Value = "Transportation"
$wsApril.Range("B10").Value = "htfg"

# Force the date-like value to stay a plain text string instead of being
# auto-converted to a date serial number.
$wsApril.Range("C10").NumberFormat = "@"
$wsApril.Range("C10").Value = "2023-04-21"
$wsApril.Range("C10").ClearFormats()

$wsApril.Range("D10").Value = 4000.0
$wsApril.Range("E10").Value = "Checkings"

# Update the monthly total (previously 5523.0) to include the new entry.
$wsApril.Range("B11").Value = 9523.0
